$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every record row.
# The whole column (rows 2-513) was bumped forward by exactly one day
# (serial date 46074 -> 46075, i.e. 2026-02-21 -> 2026-02-22) when the
# sheet was re-saved/refreshed.
$lastRow = 513

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.AddDays(1)
    }
}
